# Applies the "cryptos" price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.791.14"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.358.84"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.659"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.88%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +21.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0987"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.67%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").Value = "2.708.50"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.886"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.47%  "
$ws.Range("D17").Value = "2.354.53"
$ws.Range("E17").Value = "  +4.64%  "
$ws.Range("D18").Value = "43.664.68"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.66%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  +7.78%  "
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0702"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("E39").Value = "  +5.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.06%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.84%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0962"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("E48").Value = "  +12.43%  "
$ws.Range("D49").Value = "1.438.44"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
